$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.593000000000001
$ws.Range("B7").Value = 5.130000000000001
$ws.Range("C7").Value = -13.208
$ws.Range("C15").Value = -13.173
$ws.Range("B16").Value = 5.952
$ws.Range("D16").Value = -8.038
$ws.Range("D19").Value = -7.742
$ws.Range("C21").Value = -12.132
$ws.Range("C22").Value = -12.961
$ws.Range("C23").Value = -12.223
$ws.Range("B28").Value = 5.613
$ws.Range("B29").Value = 5.708
$ws.Range("B32").Value = 6.648000000000001
$ws.Range("C34").Value = -12.16
$ws.Range("E34").Value = 17.032
$ws.Range("D36").Value = -7.76
$ws.Range("B40").Value = 9.317
$ws.Range("C43").Value = -13.379
$ws.Range("E43").Value = 16.622
$ws.Range("C45").Value = -13.03
$ws.Range("D46").Value = -8.372
$ws.Range("E48").Value = 17.07
$ws.Range("C50").Value = -13.018
$ws.Range("D50").Value = -8.386000000000001
$ws.Range("C51").Value = -11.276
$ws.Range("B52").Value = 5.486
$ws.Range("B57").Value = 5.093000000000001
$ws.Range("B66").Value = 4.971
$ws.Range("C66").Value = -10.875
$ws.Range("C67").Value = -11.395
$ws.Range("E70").Value = 17.386
$ws.Range("E73").Value = 16.881
$ws.Range("C79").Value = -12.117
$ws.Range("C84").Value = -14.098
$ws.Range("E87").Value = 16.677
$ws.Range("C92").Value = -11.728
$ws.Range("E92").Value = 17.335
$ws.Range("D95").Value = -7.567
$ws.Range("C97").Value = -12.173
$ws.Range("D97").Value = -8.434000000000001
$ws.Range("B100").Value = 5.942
$ws.Range("E101").Value = 16.547
